$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style fix-ups first (no text/shared-string impact) ---

# D3: drop explicit style (revert to default/"Normal")
$ws.Range("D3").Style = "Normal"

# F3: drop explicit style
$ws.Range("F3").Style = "Normal"

# D4 & E4: change style from s=2 to s=1 by copying B1's format (which uses style 1)
$ws.Range("B1").Copy()
$ws.Range("D4:E4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Text / shared-string updates, in the order the new strings must be appended ---

# D3: new task text
$ws.Range("D3").Value = "Create sprites for pieces and board"

# D4: new task text
$ws.Range("D4").Value = "Represent 4x4 board and pieces via Unity"

# D5: new task text
$ws.Range("D5").Value = "Enable placing of pieces on 4x4 board via drag and drop"

# H3: contributor text (first use of "Austin/Spencer")
$ws.Range("H3").Value = "Austin/Spencer"

# I3: extra details text
$ws.Range("I3").Value = "Initial protoype and complete and functioning. Game mechanics are coming next."

# H4 / H5 reuse the same "Austin/Spencer" string
$ws.Range("H4").Value = "Austin/Spencer"
$ws.Range("H5").Value = "Austin/Spencer"

# --- Remaining numeric / cell content updates ---

# E3: fully clear (cell disappears from the row entirely)
$ws.Range("E3").Clear()

# F3: keep value 2 (style already dropped above)
$ws.Range("F3").Value = 2

# G3: keep existing style, add value 2
$ws.Range("G3").Value = 2

# F4: value 3 -> 2
$ws.Range("F4").Value = 2

# G4: add value 2 (keep existing style)
$ws.Range("G4").Value = 2

# F5: value 2 -> 3
$ws.Range("F5").Value = 3

# G5: add value 3 (keep existing style)
$ws.Range("G5").Value = 3

# --- Selection ---
$ws.Range("I3").Select()
